$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.526.44'
$ws.Range("E2").Value = '  +5.15%  '
$ws.Range("D3").Value = '2.255.14'
$ws.Range("E3").Value = '  +4.24%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.65'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.639'
$ws.Range("E6").Value = '  +2.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.23'
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.411'
$ws.Range("E9").Value = '  +3.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.28'
$ws.Range("E10").Value = '  +2.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0905'
$ws.Range("E11").Value = '  +5.49%  '
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("D13").Value = '2.587.74'
$ws.Range("E13").Value = '  +4.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.29'
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.70'
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.836'
$ws.Range("E16").Value = '  +2.74%  '
$ws.Range("E17").Value = '  +2.53%  '
$ws.Range("D18").Value = '2.252.50'
$ws.Range("E18").Value = '  +3.22%  '
$ws.Range("D19").Value = '41.444.90'
$ws.Range("E19").Value = '  +5.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.83'
$ws.Range("E20").Value = '  +2.66%  '
$ws.Range("E21").Value = '  +7.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.21'
$ws.Range("E22").Value = '  +1.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.29'
$ws.Range("E23").Value = '  +9.72%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.42'
$ws.Range("E25").Value = '  +3.44%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.40'
$ws.Range("E26").Value = '  +1.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.89'
$ws.Range("E27").Value = '  +3.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '173.51'
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  +3.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.51'
$ws.Range("E30").Value = '  +3.14%  '
$ws.Range("E31").Value = '  +2.09%  '
$ws.Range("E32").Value = '  +8.52%  '
$ws.Range("E33").Value = '  +2.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.09'
$ws.Range("E34").Value = '  +7.90%  '
$ws.Range("E35").Value = '  +3.41%  '
$ws.Range("E36").Value = '  +3.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.03'
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.86'
$ws.Range("E38").Value = '  +8.26%  '
$ws.Range("E39").Value = '  +1.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000271'
$ws.Range("E40").Value = '  +75.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.94'
$ws.Range("E42").Value = '  +14.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0242'
$ws.Range("E43").Value = '  +5.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.86'
$ws.Range("E44").Value = '  +14.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.93'
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.83'
$ws.Range("E46").Value = '  +0.93%  '
$ws.Range("E47").Value = '  +4.13%  '
$ws.Range("D48").Value = '1.513.51'
$ws.Range("E48").Value = '  -1.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0947'
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("E50").Value = '  +2.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.80'
$ws.Range("E51").Value = '  -1.02%  '
